$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.31"
$ws.Range("E2").Value = "'0.93%"
$ws.Range("D3").Value = "'29.18"
$ws.Range("E3").Value = "'2.39%"
$ws.Range("D4").Value = "'5.196"
$ws.Range("E4").Value = "'2.16%"
$ws.Range("D5").Value = "'0.06983"
$ws.Range("E5").Value = "'4.75%"
$ws.Range("D6").Value = "'7.437"
$ws.Range("E6").Value = "'1.62%"
$ws.Range("D7").Value = "'3.547"
$ws.Range("E7").Value = "'4.93%"
$ws.Range("D8").Value = "'1.406"
$ws.Range("E8").Value = "'3.29%"
$ws.Range("D9").Value = "'0.9009"
$ws.Range("E9").Value = "'-3.78%"
$ws.Range("E10").Value = "'2.59%"
$ws.Range("D11").Value = "'0.07409"
$ws.Range("E11").Value = "'17.05%"
$ws.Range("D12").Value = "'0.07732"
$ws.Range("E12").Value = "'1.44%"
$ws.Range("D13").Value = "'0.02941"
$ws.Range("E13").Value = "'2.00%"
$ws.Range("D14").Value = "'0.09018"
$ws.Range("E14").Value = "'0.68%"
$ws.Range("D15").Value = "'0.001589"
$ws.Range("E15").Value = "'-0.37%"
$ws.Range("D16").Value = "'0.0006481"
$ws.Range("E16").Value = "'1.13%"
$ws.Range("D17").Value = "'0.006152"
$ws.Range("E17").Value = "'0.75%"
$ws.Range("D18").Value = "'3.472"
$ws.Range("E18").Value = "'-0.30%"
$ws.Range("E19").Value = "'-0.11%"
$ws.Range("E20").Value = "'1.34%"
$ws.Range("D21").Value = "'0.1334"
$ws.Range("E21").Value = "'2.43%"
$ws.Range("D22").Value = "'4.021"
$ws.Range("E22").Value = "'-1.34%"
$ws.Range("D24").Value = "'0.04514"
$ws.Range("E24").Value = "'1.11%"
$ws.Range("D25").Value = "'0.001207"
$ws.Range("E25").Value = "'2.34%"
$ws.Range("D26").Value = "'0.004243"
$ws.Range("E26").Value = "'-5.02%"
$ws.Range("E27").Value = "'-6.26%"
$ws.Range("D28").Value = "'0.0001665"
$ws.Range("E28").Value = "'3.31%"
$ws.Range("D40").Value = "'0.04355"
$ws.Range("E40").Value = "'4.59%"
$ws.Range("D41").Value = "'0.006931"
$ws.Range("E41").Value = "'2.57%"
$ws.Range("D42").Value = "'0.1245"
$ws.Range("E42").Value = "'-0.06%"
$ws.Range("D43").Value = "'0.002065"
$ws.Range("E43").Value = "'2.62%"
$ws.Range("D44").Value = "'0.01156"
$ws.Range("E44").Value = "'0.36%"
$ws.Range("D45").Value = "'0.00005815"
$ws.Range("E45").Value = "'2.80%"
$ws.Range("D47").Value = "'0.01304"
$ws.Range("E47").Value = "'0.15%"
